$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-42) holds the "Förändrad" (Changed) date.
# Update its serial value from 45715 (2025-02-27) to 45716 (2025-02-28).
for ($row = 2; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45715) {
        $cell.Value2 = 45716
    }
}
